$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 909414.0600000001
$ws.Range("I2").Value = 1000277.8
$ws.Range("J2").Value = 777
$ws.Range("K2").Value = 1000277.8
$ws.Range("L2").Value = 777
$ws.Range("M2").Value = -1000164.8
$ws.Range("N2").Value = -1003

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2881.9375
$ws.Range("I38").Value = 411.4
$ws.Range("J38").Value = 6999.5
$ws.Range("K38").Value = 1234.2
$ws.Range("L38").Value = 20998.5
$ws.Range("M38").Value = -862.1999999999998
$ws.Range("N38").Value = -21742.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1998.6666
$ws.Range("J39").Value = 1998.6666
$ws.Range("L39").Value = 5995.9998
$ws.Range("N39").Value = -6587.9998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3180.8635
$ws.Range("I51").Value = 2500
$ws.Range("K51").Value = 2500
$ws.Range("M51").Value = -2016

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1824.7142
$ws.Range("J129").Value = 2330.6667
$ws.Range("L129").Value = 6992.000100000001
$ws.Range("N129").Value = -16992.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3050.1875
$ws.Range("I132").Value = 2407.923
$ws.Range("K132").Value = 7223.768999999999
$ws.Range("M132").Value = -4693.768999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 74999.5
$ws.Range("J133").Value = 74999.5
$ws.Range("L133").Value = 74999.5
$ws.Range("N133").Value = -85119.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 8020.364
$ws.Range("I141").Value = 5181.125
$ws.Range("K141").Value = 15543.375
$ws.Range("M141").Value = -10363.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2808.7
$ws.Range("I2").Value = 2242.8948
$ws.Range("J2").Value = 3786
$ws.Range("K2").Value = 2242.8948
$ws.Range("L2").Value = 3786
$ws.Range("M2").Value = -2129.8948
$ws.Range("N2").Value = -4012

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 87.40000000000001
$ws.Range("I5").Value = 87.40000000000001
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 87.40000000000001
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 24.59999999999999
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2831.12
$ws.Range("I32").Value = 2707.9348
$ws.Range("K32").Value = 2707.9348
$ws.Range("M32").Value = -2420.9348

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 3734.6667
$ws.Range("I104").Value = 3734.6667
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 3734.6667
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -240.6667000000002
$ws.Range("N104").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2808.7
$ws.Range("I116").Value = 2242.8948
$ws.Range("J116").Value = 3786
$ws.Range("K116").Value = 2242.8948
$ws.Range("L116").Value = 3786
$ws.Range("M116").Value = 51.10519999999997
$ws.Range("N116").Value = -8374

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2948.625
$ws.Range("I122").Value = 2948.625
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8845.875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6395.875
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2808.7
$ws.Range("I3").Value = 2242.8948
$ws.Range("J3").Value = 3786
$ws.Range("K3").Value = 2242.8948
$ws.Range("L3").Value = 3786
$ws.Range("M3").Value = -2128.8948
$ws.Range("N3").Value = -4014

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 87.40000000000001
$ws.Range("I4").Value = 87.40000000000001
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 87.40000000000001
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 27.59999999999999
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1664.0769
$ws.Range("I20").Value = 1828.3334
$ws.Range("J20").Value = 1294.5
$ws.Range("K20").Value = 1828.3334
$ws.Range("L20").Value = 1294.5
$ws.Range("M20").Value = -1581.3334
$ws.Range("N20").Value = -1788.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 21126.125
$ws.Range("J95").Value = 21126.125
$ws.Range("L95").Value = 21126.125
$ws.Range("N95").Value = -26618.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 540.3570999999999
$ws.Range("I5").Value = 195.77777
$ws.Range("J5").Value = 1160.6
$ws.Range("K5").Value = 195.77777
$ws.Range("L5").Value = 1160.6
$ws.Range("M5").Value = -83.77777
$ws.Range("N5").Value = -1384.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 2039.8
$ws.Range("J12").Value = 2499.75
$ws.Range("L12").Value = 2499.75
$ws.Range("N12").Value = -2839.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4571.385
$ws.Range("I31").Value = 1805
$ws.Range("K31").Value = 1805
$ws.Range("M31").Value = -1510

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4571.385
$ws.Range("I34").Value = 1805
$ws.Range("K34").Value = 1805
$ws.Range("M34").Value = -1603

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5725.5
$ws.Range("I86").Value = 5725.5
$ws.Range("K86").Value = 5725.5
$ws.Range("M86").Value = -4602.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 5725.5
$ws.Range("I89").Value = 5725.5
$ws.Range("K89").Value = 28627.5
$ws.Range("M89").Value = -23011.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H114").Value = 18066.555
$ws.Range("I114").Value = 9665.833000000001
$ws.Range("J114").Value = 19358.975
$ws.Range("K114").Value = 9665.833000000001
$ws.Range("L114").Value = 19358.975
$ws.Range("M114").Value = -5326.833000000001
$ws.Range("N114").Value = -28036.975

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 60000
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1839.72
$ws.Range("I132").Value = 1823.2858
$ws.Range("K132").Value = 5469.857400000001
$ws.Range("M132").Value = -2939.857400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 250000
$ws.Range("I140").Value = 250000
$ws.Range("K140").Value = 250000
$ws.Range("M140").Value = -244820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 22.8
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 33.333332
$ws.Range("K2").Value = 42
$ws.Range("L2").Value = 199.999992
$ws.Range("M2").Value = 71
$ws.Range("N2").Value = -425.999992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1867524.1
$ws.Range("J4").Value = 25750000
$ws.Range("L4").Value = 77250000
$ws.Range("N4").Value = -77250224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1637.4667
$ws.Range("I5").Value = 946.3333
$ws.Range("J5").Value = 2098.2222
$ws.Range("K5").Value = 2838.9999
$ws.Range("L5").Value = 6294.6666
$ws.Range("M5").Value = -2726.9999
$ws.Range("N5").Value = -6518.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3911.6365
$ws.Range("I34").Value = 283.33334
$ws.Range("J34").Value = 4274.467
$ws.Range("K34").Value = 850.0000200000001
$ws.Range("L34").Value = 12823.401
$ws.Range("M34").Value = -766.0000200000001
$ws.Range("N34").Value = -12991.401

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 12705.883
$ws.Range("I88").Value = 4400
$ws.Range("J88").Value = 16166.667
$ws.Range("K88").Value = 13200
$ws.Range("L88").Value = 48500.001
$ws.Range("M88").Value = -12772
$ws.Range("N88").Value = -49356.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 12705.883
$ws.Range("I91").Value = 4400
$ws.Range("J91").Value = 16166.667
$ws.Range("K91").Value = 13200
$ws.Range("L91").Value = 48500.001
$ws.Range("M91").Value = -11718
$ws.Range("N91").Value = -51464.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2853.625
$ws.Range("I122").Value = 1010.1429
$ws.Range("J122").Value = 4287.4443
$ws.Range("K122").Value = 9091.286100000001
$ws.Range("L122").Value = 38586.9987
$ws.Range("M122").Value = -6641.286100000001
$ws.Range("N122").Value = -43486.9987

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1637.4667
$ws.Range("I135").Value = 946.3333
$ws.Range("J135").Value = 2098.2222
$ws.Range("K135").Value = 8516.9997
$ws.Range("L135").Value = 18883.9998
$ws.Range("M135").Value = -5981.9997
$ws.Range("N135").Value = -23953.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 3005
$ws.Range("I140").Value = 2079.8
$ws.Range("J140").Value = 4547
$ws.Range("K140").Value = 6239.400000000001
$ws.Range("L140").Value = 13641
$ws.Range("M140").Value = -1059.400000000001
$ws.Range("N140").Value = -24001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4515.8335
$ws.Range("I80").Value = 2398.3333
$ws.Range("J80").Value = 4939.3335
$ws.Range("K80").Value = 2398.3333
$ws.Range("L80").Value = 4939.3335
$ws.Range("M80").Value = -1400.3333
$ws.Range("N80").Value = -6935.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4515.8335
$ws.Range("I83").Value = 2398.3333
$ws.Range("J83").Value = 4939.3335
$ws.Range("K83").Value = 11991.6665
$ws.Range("L83").Value = 24696.6675
$ws.Range("M83").Value = -6999.666499999999
$ws.Range("N83").Value = -34680.6675

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 10583.333
$ws.Range("J92").Value = 10583.333
$ws.Range("L92").Value = 10583.333
$ws.Range("N92").Value = -14327.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2238.6072
$ws.Range("I102").Value = 2287.7307
$ws.Range("J102").Value = 1600
$ws.Range("K102").Value = 2287.7307
$ws.Range("L102").Value = 1600
$ws.Range("M102").Value = -665.7307000000001
$ws.Range("N102").Value = -4844

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2170.1052
$ws.Range("I132").Value = 2284.3125
$ws.Range("K132").Value = 6852.9375
$ws.Range("M132").Value = -4322.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 40000
$ws.Range("J134").Value = 40000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -125070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 39698.875
$ws.Range("J136").Value = 39698.875
$ws.Range("L136").Value = 119096.625
$ws.Range("N136").Value = -124196.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 10375
$ws.Range("I13").Value = 6250
$ws.Range("J13").Value = 14500
$ws.Range("K13").Value = 6250
$ws.Range("L13").Value = 14500
$ws.Range("M13").Value = -6110
$ws.Range("N13").Value = -14780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2762.7727
$ws.Range("I22").Value = 2084.1
$ws.Range("J22").Value = 3328.3333
$ws.Range("K22").Value = 2084.1
$ws.Range("L22").Value = 3328.3333
$ws.Range("M22").Value = -1789.1
$ws.Range("N22").Value = -3918.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2762.7727
$ws.Range("I27").Value = 2084.1
$ws.Range("J27").Value = 3328.3333
$ws.Range("K27").Value = 2084.1
$ws.Range("L27").Value = 3328.3333
$ws.Range("M27").Value = -1977.1
$ws.Range("N27").Value = -3542.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1214.826
$ws.Range("J55").Value = 2199.4
$ws.Range("L55").Value = 2199.4
$ws.Range("N55").Value = -2545.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 9890.666999999999
$ws.Range("I68").Value = 2379.4
$ws.Range("J68").Value = 15255.857
$ws.Range("K68").Value = 2379.4
$ws.Range("L68").Value = 15255.857
$ws.Range("M68").Value = -1630.4
$ws.Range("N68").Value = -16753.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 9890.666999999999
$ws.Range("I71").Value = 2379.4
$ws.Range("J71").Value = 15255.857
$ws.Range("K71").Value = 11897
$ws.Range("L71").Value = 76279.285
$ws.Range("M71").Value = -8153
$ws.Range("N71").Value = -83767.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4233.0586
$ws.Range("I122").Value = 3696.9285
$ws.Range("J122").Value = 6735
$ws.Range("K122").Value = 11090.7855
$ws.Range("L122").Value = 20205
$ws.Range("M122").Value = -8640.7855
$ws.Range("N122").Value = -25105

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1408.3846
$ws.Range("I13").Value = 215
$ws.Range("J13").Value = 4093.5
$ws.Range("K13").Value = 215
$ws.Range("L13").Value = 4093.5
$ws.Range("M13").Value = -75
$ws.Range("N13").Value = -4373.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 39998.5
$ws.Range("J70").Value = 39998.5
$ws.Range("L70").Value = 39998.5
$ws.Range("N70").Value = -40628.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 39998.5
$ws.Range("J73").Value = 39998.5
$ws.Range("L73").Value = 39998.5
$ws.Range("N73").Value = -42182.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6247.5713
$ws.Range("I126").Value = 7150.6
$ws.Range("J126").Value = 3990
$ws.Range("K126").Value = 21451.8
$ws.Range("L126").Value = 11970
$ws.Range("M126").Value = -18981.8
$ws.Range("N126").Value = -16910
